# edit.ps1 - applies the Fig1.pptx commit:
#   1) datetimeFigureOut date placeholders: "23.10.2023" -> "01.11.2023"
#      (slide master, all slide layouts, notes master)
#   2) slide 1 text box: drop the "Rb+" run and fix " Cs" -> "Cs"

$p = $ppt.ActivePresentation

$oldDate = "23.10.2023"
$newDate = "01.11.2023"

function Update-DatePlaceholder($shp) {
    if ($shp.HasTextFrame -and $shp.Type -eq 14) {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# --- Slide Master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DatePlaceholder $master.Shapes.Item($i)
}

# --- All Slide Layouts ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DatePlaceholder $layout.Shapes.Item($i)
    }
}

# --- Notes Master ---
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    Update-DatePlaceholder $notesMaster.Shapes.Item($i)
}

# --- Any date placeholders on the slides themselves ---
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        Update-DatePlaceholder $slide.Shapes.Item($i)
    }
}

# --- Slide 1: ion textbox, remove "Rb+" and tighten " Cs" -> "Cs" ---
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "*Rb+ Cs+*") {
            $full = $tr.Text
            $idx = $full.IndexOf("Rb+ Cs")
            # Characters() is 1-based
            $rbStart = $idx + 1
            # Remove "Rb+" (3 characters)
            $tr.Characters($rbStart, 3).Text = ""
            # Now remove the leftover leading space before "Cs" (was " Cs")
            $tr.Characters($rbStart, 1).Text = ""
        }
    }
}
